$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Updated loading_percent values for rows 2-25 (data index 0-23),
# columns B,C,E,F,G,H,I,L,M  -- "case with 380 kV done"
# Flat array: 24 rows x 9 columns, row-major order.
$values = @(
    16.06204481467558, 9.812188730489245, 11.5997085747307, 16.86991607391233, 35.35228157417627, 15.89065115055791, 24.39842647788738, 10.11955124519447, 15.17540566081926,  # row 2 (A2=0)
    15.54408127764977, 9.436218910696313, 11.62617660362494, 15.89584955866808, 35.48241100409955, 15.96987771850882, 24.56573891250344, 10.13051417286695, 15.06991966843087,  # row 3 (A3=1)
    15.21934594939388, 9.19578111948049, 11.64334190461561, 15.26997757108491, 35.57943185957705, 16.02255113301347, 24.67524438378815, 10.13874153750899, 15.00713458932451,  # row 4 (A4=2)
    15.08553050268584, 9.095475666626919, 11.65056735504728, 15.00819731993403, 35.62323429148654, 16.0450257032762, 24.72156874038151, 10.1424706361239, 14.98206857425658,  # row 5 (A5=3)
    15.06322720577416, 9.078682341868173, 11.65178107376305, 14.96433081551593, 35.63076408514712, 16.04881848207676, 24.72936343647573, 10.1431125865105, 14.977938338527,  # row 6 (A6=4)
    15.21754699327366, 9.194437656698247, 11.64343841559326, 15.26647399323137, 35.58000537715773, 16.02285014876699, 24.675862251827, 10.13879030530096, 15.00679441061683,  # row 7 (A7=5)
    15.88495878377198, 9.684593467095569, 11.60864548591103, 16.53996406344768, 35.3935737254694, 15.91713057050847, 24.45470744349881, 10.12302087532827, 15.13863494965816,  # row 8 (A8=6)
    17.13249556972094, 10.56642648913219, 11.54763732432126, 19.00274580682531, 35.16547265956542, 15.7419160958752, 24.07494638829064, 10.10395837388312, 15.41197557646064,  # row 9 (A9=7)
    18.00207716687726, 11.16224439437437, 11.50717453844395, 20.67494806633232, 35.08375903711858, 15.63294982712086, 23.82904395906861, 10.09716899829199, 15.6205234149592,  # row 10 (A10=8)
    18.38575351758728, 11.42138521144217, 11.4897047750375, 21.3917225636224, 35.0656063066787, 15.58771194037852, 23.72442021807237, 10.09564257166035, 15.71678767825664,  # row 11 (A11=9)
    18.52920770784186, 11.51776282472831, 11.48322348964036, 21.65686569030329, 35.06149278770916, 15.571207958805, 23.68584787975128, 10.09528863814354, 15.75341634089101,  # row 12 (A12=10)
    18.49839575701507, 11.49708476843154, 11.48461339421146, 21.60004134736742, 35.06225556451617, 15.57473445698415, 23.69410849374411, 10.0953549069366, 15.74552030014693,  # row 13 (A13=11)
    18.39759303328648, 11.42934961538199, 11.48916887093473, 21.4136618050453, 35.06521242655311, 15.58634156773002, 23.72122584475524, 10.0956089655727, 15.71979776796791,  # row 14 (A14=12)
    18.33560593816621, 11.38763039158548, 11.49197668260456, 21.29868154950795, 35.06738378820219, 15.59353297544491, 23.73797244511304, 10.09579374982312, 15.70406407841976,  # row 15 (A15=13)
    17.97675258299385, 11.14506572156003, 11.50833503128493, 20.62722412089977, 35.0853303469277, 15.63599371058905, 23.8360274900176, 10.09730014319013, 15.61425824528412,  # row 16 (A16=14)
    17.75346655640784, 10.99318021306884, 11.51860989851871, 20.20408069597325, 35.10123061928323, 15.66315405371853, 23.898038891199, 10.09862400047147, 15.55950473035677,  # row 17 (A17=15)
    17.62392469172019, 10.90470170935451, 11.52460795704696, 19.95656407809801, 35.11216490903217, 15.6791834310677, 23.93438718773503, 10.09953251699292, 15.52814476050429,  # row 18 (A18=16)
    17.57987679682372, 10.87455378843952, 11.52665396795667, 19.87204792380568, 35.11617347296246, 15.6846805709768, 23.94681089091762, 10.09986539932053, 15.51755038402858,  # row 19 (A19=17)
    17.77735193967665, 11.00946471748925, 11.51750699459461, 20.24955283636154, 35.09935266930982, 15.66022059027798, 23.89136715364866, 10.09846785616022, 15.56531977195759,  # row 20 (A20=18)
    18.42725197934266, 11.44929296308559, 11.48782718135877, 21.46857628470577, 35.06426881905356, 15.58291524138694, 23.71323237732351, 10.09552826529435, 15.72734853772396,  # row 21 (A21=19)
    18.84125127573572, 11.72651082244485, 11.46921127437288, 22.22866616901552, 35.05743917197092, 15.53604593526572, 23.6029129531657, 10.09491294182625, 15.83425476448184,  # row 22 (A22=20)
    18.62131300953667, 11.57950329695283, 11.47907561730085, 21.82633154458858, 35.05960367361225, 15.5607253551199, 23.66123229700326, 10.09512206230238, 15.77711271854736,  # row 23 (A23=21)
    17.76655699804536, 11.00210609526142, 11.51800533427503, 20.22900810905287, 35.10019610723755, 15.66154551806059, 23.89438127501711, 10.09853798978136, 15.56269042042772,  # row 24 (A24=22)
    16.80262022683755, 10.33677843651951, 11.56337298888685, 18.34778573295695, 35.21223690985264, 15.78585965406964, 24.17188547871214, 10.10784671996139, 15.33657696960899,  # row 25 (A25=23)
)

$cols = @(2,3,5,6,7,8,9,12,13)   # B,C,E,F,G,H,I,L,M

$idx = 0
for ($r = 2; $r -le 25; $r++) {
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Cells.Item($r, $cols[$j]).Value = $values[$idx]
        $idx = $idx + 1
    }
}
